$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers (columns shift right by one, new F1 header added) ---
# A1 needs the same header formatting (bold/border/alignment) as the rest
# of row 1 -- clone it from an existing formatted header cell so it reuses
# the same style record instead of creating a new, merely-equivalent one.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("A1").Value = "groups"
$ws.Range("B1").Value = "independant_variables"
$ws.Range("C1").Value = "paired"
$ws.Range("D1").Value = "parametric"
$ws.Range("E1").Value = "label"
$ws.Range("F1").Value = "control_group_id"

# --- Row 2: data values shift columns too ---
$ws.Range("A2").Value = "1, 5, 3, 4"
$ws.Range("B2").Value = "TCB2, MDL"
$ws.Range("C2").Value = $false
$ws.Range("D2").Value = $true
$ws.Range("E2").Value = "agonist antagonist"
$ws.Range("F2").Value = 1

# A2 previously carried the header-style formatting (s="1"); the new
# layout has no special formatting on row 2, so clear it.
$ws.Range("A2").ClearFormats()
